$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# New values for rows 2-8, columns A, B, E, F, G, H, Q, R
# (row 7 is unchanged and omitted)
$data = @{
    2 = @{ A = 112243463; B = 77403;  E = 228912; F = "Mörk kolflarnlav"; G = "Carbonicola myrmecina"; H = "(Ach.) Bendiksby & Timdal"; Q = 410608; R = 6710914 }
    3 = @{ A = 112243468; B = 77650;  E = 6425;   F = "Garnlav";          G = "Alectoria sarmentosa";   H = "(Ach.) Ach.";                 Q = 410566; R = 6710872 }
    4 = @{ A = 112243462; B = 78242;  E = 6453;   F = "Vedskivlav";       G = "Hertelidea botryosa";    H = "(Fr.) Printzen & Kantvilas";  Q = 410608; R = 6710914 }
    5 = @{ A = 112243461; B = 77403;  E = 228912; F = "Mörk kolflarnlav"; G = "Carbonicola myrmecina"; H = "(Ach.) Bendiksby & Timdal"; Q = 410598; R = 6710899 }
    6 = @{ A = 112243460; B = 78242;  E = 6453;   F = "Vedskivlav";       G = "Hertelidea botryosa";    H = "(Fr.) Printzen & Kantvilas";  Q = 410598; R = 6710899 }
    8 = @{ A = 112243469; B = 77650;  E = 6425;   F = "Garnlav";          G = "Alectoria sarmentosa";   H = "(Ach.) Ach.";                 Q = 410486; R = 6710828 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
}
